# Add team record (Wins/Losses/Ties) columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from the last existing header cell (AC1) onto
# the three new header cells so they match the other headers (bold, border,
# centered).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record for every data row (rows 2-54).
for ($r = 2; $r -le 54; $r++) {
    $ws.Cells.Item($r, 30).Value = 80   # column AD - Wins
    $ws.Cells.Item($r, 31).Value = 82   # column AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # column AF - Ties
}
